# Add data for 2022-05-26: update the "through" date from 05-17 to 05-18
# and bump the May / Total figures for a handful of years.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet/tab
$ws.Name = "Through 2022-05-18"

# 2. Update the "May (through 05-17)" label in column A
$ws.Range("A6").Value = "May (through 05-18)"

# 3. Update the May row (row 6) figures
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 37
$ws.Range("F6").Value = 26
$ws.Range("G6").Value = 35
$ws.Range("H6").Value = 66

# 4. Update the Total row (row 7) figures
$ws.Range("C7").Value = 191
$ws.Range("D7").Value = 290
$ws.Range("F7").Value = 181
$ws.Range("G7").Value = 297
$ws.Range("H7").Value = 589
